# Natmi following Dr Hou advice
# Rewrite the LR-pairs data: add the "sCs" cluster as a new source/target
# cluster (3 clusters x 3 clusters = 9 rows instead of 2 clusters x ... rows),
# and update all the numeric values to match the recomputed statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row=2;  A="ECs";  B="Angpt1"; C="Tek"; D="ECs";  E=1; F=0.3333333333333333; G=0.110507; H=0.331521; I=0.004605687348208628; J=0.004605687348208628; K=3; L=1; M=64.221457; N=192.664371; O=0.8366610301096816; P=0.8366610301096816; Q=7.096920548699001; R=63.87228493829101; S=0.003853399121115358; T=0.003853399121115358 },
    @{ Row=3;  A="ECs";  B="Angpt1"; C="Tek"; D="FAPs"; E=1; F=0.3333333333333333; G=0.110507; H=0.331521; I=0.004605687348208628; J=0.004605687348208628; K=3; L=1; M=10.582537; N=31.747611; O=0.1378666371208897; P=0.1378666371208896; Q=1.169444416259; R=10.524999746331; S=0.0006349706263277515; T=0.0006349706263277514 },
    @{ Row=4;  A="ECs";  B="Angpt1"; C="Tek"; D="sCs";  E=1; F=0.3333333333333333; G=0.110507; H=0.331521; I=0.004605687348208628; J=0.004605687348208628; K=3; L=1; M=1.955236666666667; N=5.86571; O=0.02547233276942866; P=0.02547233276942866; Q=0.2160673383233334; R=1.94460604491; S=0.0001173176007655176; T=0.0001173176007655176 },
    @{ Row=5;  A="FAPs"; B="Angpt1"; C="Tek"; D="ECs";  E=3; F=1; G=22.14783133333333; H=66.443494; I=0.9230726249214253; J=0.9230726249214253; K=3; L=1; M=64.221457; N=192.664371; O=0.8366610301096816; P=0.8366610301096816; Q=1422.365997616919; R=12801.29397855228; S=0.7722988932328074; T=0.7722988932328074 },
    @{ Row=6;  A="FAPs"; B="Angpt1"; C="Tek"; D="FAPs"; E=3; F=1; G=22.14783133333333; H=66.443494; I=0.9230726249214253; J=0.9230726249214253; K=3; L=1; M=10.582537; N=31.747611; O=0.1378666371208897; P=0.1378666371208896; Q=234.3802445547593; R=2109.422200992834; S=0.1272609186162692; T=0.1272609186162692 },
    @{ Row=7;  A="FAPs"; B="Angpt1"; C="Tek"; D="sCs";  E=3; F=1; G=22.14783133333333; H=66.443494; I=0.9230726249214253; J=0.9230726249214253; K=3; L=1; M=1.955236666666667; N=5.86571; O=0.02547233276942866; P=0.02547233276942866; Q=43.30425191008222; R=389.73826719074; S=0.02351281307234855; T=0.02351281307234855 },
    @{ Row=8;  A="sCs";  B="Angpt1"; C="Tek"; D="ECs";  E=3; F=1; G=1.735257333333333; H=5.205772; I=0.07232168773036617; J=0.07232168773036617; K=3; L=1; M=64.221457; N=192.664371; O=0.8366610301096816; P=0.8366610301096816; Q=111.4407542166013; R=1002.966787949412; S=0.06050873775575888; T=0.06050873775575888 },
    @{ Row=9;  A="sCs";  B="Angpt1"; C="Tek"; D="FAPs"; E=3; F=1; G=1.735257333333333; H=5.205772; I=0.07232168773036617; J=0.07232168773036617; K=3; L=1; M=10.582537; N=31.747611; O=0.1378666371208897; P=0.1378666371208896; Q=18.36342493452133; R=165.270824410692; S=0.009970747878292691; T=0.009970747878292689 },
    @{ Row=10; A="sCs";  B="Angpt1"; C="Tek"; D="sCs";  E=3; F=1; G=1.735257333333333; H=5.205772; I=0.07232168773036617; J=0.07232168773036617; K=3; L=1; M=1.955236666666667; N=5.86571; O=0.02547233276942866; P=0.02547233276942866; Q=3.392838764235556; R=30.53554887812; S=0.001842202096314593; T=0.001842202096314592 }
)

foreach ($rec in $data) {
    $r = $rec.Row
    $ws.Range("A$r").Value = $rec.A
    $ws.Range("B$r").Value = $rec.B
    $ws.Range("C$r").Value = $rec.C
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("E$r").Value = $rec.E
    $ws.Range("F$r").Value = $rec.F
    $ws.Range("G$r").Value = $rec.G
    $ws.Range("H$r").Value = $rec.H
    $ws.Range("I$r").Value = $rec.I
    $ws.Range("J$r").Value = $rec.J
    $ws.Range("K$r").Value = $rec.K
    $ws.Range("L$r").Value = $rec.L
    $ws.Range("M$r").Value = $rec.M
    $ws.Range("N$r").Value = $rec.N
    $ws.Range("O$r").Value = $rec.O
    $ws.Range("P$r").Value = $rec.P
    $ws.Range("Q$r").Value = $rec.Q
    $ws.Range("R$r").Value = $rec.R
    $ws.Range("S$r").Value = $rec.S
    $ws.Range("T$r").Value = $rec.T
}
